$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are stored as text, preserving exact
# formatting (e.g. trailing zeros, double-dot thousand separators) just
# like the source data, instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "24.410.87"
$ws.Cells.Item(2, 5).Value = "  +8.71%  "

$ws.Cells.Item(3, 4).Value = "1.678.12"
$ws.Cells.Item(3, 5).Value = "  +4.58%  "

$ws.Cells.Item(4, 4).Value = "1.007"
$ws.Cells.Item(4, 5).Value = "  +0.33%  "

$ws.Cells.Item(5, 4).Value = "306.62"
$ws.Cells.Item(5, 5).Value = "  +0.31%  "

$ws.Cells.Item(6, 4).Value = "0.9997"
$ws.Cells.Item(6, 5).Value = "  +0.77%  "

$ws.Cells.Item(7, 4).Value = "0.3702"
$ws.Cells.Item(7, 5).Value = "  +0.19%  "

$ws.Cells.Item(8, 4).Value = "0.3447"
$ws.Cells.Item(8, 5).Value = "  +1.23%  "

$ws.Cells.Item(9, 4).Value = "48.30"
$ws.Cells.Item(9, 5).Value = "  +13.89%  "

$ws.Cells.Item(10, 4).Value = "1.179"
$ws.Cells.Item(10, 5).Value = "  +3.22%  "

$ws.Cells.Item(11, 4).Value = "0.07266"
$ws.Cells.Item(11, 5).Value = "  +2.59%  "

$ws.Cells.Item(12, 4).Value = "1.001"
$ws.Cells.Item(12, 5).Value = "  +0.16%  "

$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(13, 4).Value = "6.140"
$ws.Cells.Item(13, 5).Value = "  +3.06%  "

$ws.Cells.Item(14, 2).Value = "Solana"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(14, 4).Value = "20.35"
$ws.Cells.Item(14, 5).Value = "  +2.72%  "

$ws.Cells.Item(15, 4).Value = "6.745"
$ws.Cells.Item(15, 5).Value = "  +1.30%  "

$ws.Cells.Item(16, 4).Value = "1.678.30"
$ws.Cells.Item(16, 5).Value = "  +4.64%  "

$ws.Cells.Item(17, 4).Value = "0.00001107"
$ws.Cells.Item(17, 5).Value = "  +1.23%  "

$ws.Cells.Item(18, 4).Value = "0.9992"
$ws.Cells.Item(18, 5).Value = "  +0.65%  "

$ws.Cells.Item(19, 4).Value = "0.06722"
$ws.Cells.Item(19, 5).Value = "  -1.52%  "

$ws.Cells.Item(20, 4).Value = "81.09"
$ws.Cells.Item(20, 5).Value = "  +3.76%  "

$ws.Cells.Item(21, 4).Value = "16.44"
$ws.Cells.Item(21, 5).Value = "  +1.67%  "

$ws.Cells.Item(22, 4).Value = "6.099"
$ws.Cells.Item(22, 5).Value = "  +0.74%  "

$ws.Cells.Item(23, 4).Value = "12.02"
$ws.Cells.Item(23, 5).Value = "  +1.18%  "

$ws.Cells.Item(24, 4).Value = "24.382.23"
$ws.Cells.Item(24, 5).Value = "  +8.67%  "

$ws.Cells.Item(25, 4).Value = "2.438"
$ws.Cells.Item(25, 5).Value = "  +1.18%  "

$ws.Cells.Item(26, 4).Value = "2.675"
$ws.Cells.Item(26, 5).Value = "  +4.90%  "

$ws.Cells.Item(27, 4).Value = "152.38"
$ws.Cells.Item(27, 5).Value = "  +0.71%  "

$ws.Cells.Item(28, 4).Value = "19.48"
$ws.Cells.Item(28, 5).Value = "  -0.77%  "

$ws.Cells.Item(29, 4).Value = "1.861.88"
$ws.Cells.Item(29, 5).Value = "  +4.44%  "

$ws.Cells.Item(30, 4).Value = "126.80"
$ws.Cells.Item(30, 5).Value = "  +4.67%  "

$ws.Cells.Item(31, 4).Value = "6.357"
$ws.Cells.Item(31, 5).Value = "  +3.06%  "

$ws.Cells.Item(32, 4).Value = "4.042"
$ws.Cells.Item(32, 5).Value = "  -3.68%  "

$ws.Cells.Item(33, 4).Value = "0.9807"
$ws.Cells.Item(33, 5).Value = "  +2.50%  "

$ws.Cells.Item(34, 4).Value = "0.08437"
$ws.Cells.Item(34, 5).Value = "  +1.85%  "

$ws.Cells.Item(35, 4).Value = "1.700"
$ws.Cells.Item(35, 5).Value = "  +3.75%  "

$ws.Cells.Item(36, 4).Value = "12.59"
$ws.Cells.Item(36, 5).Value = "  +4.84%  "

$ws.Cells.Item(37, 4).Value = "0.06522"
$ws.Cells.Item(37, 5).Value = "  +6.23%  "

$ws.Cells.Item(38, 4).Value = "5.359"
$ws.Cells.Item(38, 5).Value = "  +1.00%  "

$ws.Cells.Item(39, 4).Value = "8.907"
$ws.Cells.Item(39, 5).Value = "  +3.03%  "

$ws.Cells.Item(40, 5).Value = "  +4.05%  "

$ws.Cells.Item(41, 4).Value = "1.261"
$ws.Cells.Item(41, 5).Value = "  -0.37%  "

$ws.Cells.Item(42, 4).Value = "0.2113"
$ws.Cells.Item(42, 5).Value = "  +3.86%  "

$ws.Cells.Item(43, 4).Value = "0.6178"
$ws.Cells.Item(43, 5).Value = "  +4.00%  "

$ws.Cells.Item(44, 4).Value = "0.9989"
$ws.Cells.Item(44, 5).Value = "  +0.72%  "

$ws.Cells.Item(45, 4).Value = "13.23"
$ws.Cells.Item(45, 5).Value = "  +0.74%  "

$ws.Cells.Item(46, 2).Value = "PancakeSwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(46, 4).Value = "3.760"
$ws.Cells.Item(46, 5).Value = "  -2.23%  "

$ws.Cells.Item(47, 2).Value = "Decentraland"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(47, 4).Value = "0.5953"
$ws.Cells.Item(47, 5).Value = "  +4.14%  "

$ws.Cells.Item(48, 4).Value = "127.36"
$ws.Cells.Item(48, 5).Value = "  -0.40%  "

$ws.Cells.Item(49, 4).Value = "2.029"
$ws.Cells.Item(49, 5).Value = "  +2.06%  "

$ws.Cells.Item(50, 4).Value = "0.07203"
$ws.Cells.Item(50, 5).Value = "  +5.72%  "

$ws.Cells.Item(51, 4).Value = "75.87"
$ws.Cells.Item(51, 5).Value = "  +2.31%  "
